$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

# --- Update source data -------------------------------------------------
# "SP Reales completados" for the first two members now match their
# estimates, so the deviation drops to 0 and every dependent formula
# (Rendimiento / Desviacion in Table2/Table3/Table35, plus the chart
# series that reads from them) recalculates automatically.
$ws.Range("C2").Value = 10.5
$ws.Range("C3").Value = 12

# --- Chart: value-axis title -------------------------------------------
# Shorten the axis title and make it horizontal instead of rotated.
$co = $ws.ChartObjects(1)
$chart = $co.Chart
$valueAxis = $chart.Axes(2)
$axisTitle = $valueAxis.AxisTitle
$axisTitle.Text = "%"
$axisTitle.Orientation = 0

# --- Selection ------------------------------------------------------------
$ws.Range("C5").Select()
